$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop the now-unused trailing columns (F:N). The new table only
#    needs 2017-2021 (cols D:H) instead of 2008-2018 (cols D:N).
# ---------------------------------------------------------------------
$ws.Range("F1:N8").EntireColumn.Delete()

# ---------------------------------------------------------------------
# 2. Title row (row 1) - drop the trailing "*" from the heading and
#    grow the row height a touch to fit the rewrapped text.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "9.5.1 ИДП га болгон тажрыйбалык-конструктордук жумуштун жана илимий изилдөөнүн чыгымдарынын үлүшү"
$ws.Range("B1").Value = "9.5.1  Доля расходов на научно-исследовательские и опытно-конструкторские работы в ВВП"
$ws.Range("C1").Value = "9.5.1 Research and development expenditure as a proportion of GDP"
$ws.Rows.Item(1).RowHeight = 43.5

# ---------------------------------------------------------------------
# 3. Header row (row 4) - recent years + a new "2021*" column that
#    needs a dedicated right-aligned bold style.
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 2017
$ws.Range("E4").Value = 2018

$ws.Range("E4").Copy()
$ws.Range("F4:H4").PasteSpecial(-4122)
$ws.Range("F4").Value = 2019
$ws.Range("G4").Value = 2020
$ws.Range("H4").Value = "2021*"
$ws.Range("H4").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 4. Data row (row 5) - new figures for 2017-2021, and extend the
#    bottom-border style that row 3/5 carry out through column H.
# ---------------------------------------------------------------------
$ws.Range("D5").Value = 0.11
$ws.Range("E5").Value = 0.1

$ws.Range("E5").Copy()
$ws.Range("F5:H5").PasteSpecial(-4122)
$ws.Range("F5").Value = 0.09
$ws.Range("G5").Value = 0.09
$ws.Range("H5").Value = 0.08

$ws.Range("E3").Copy()
$ws.Range("F3:H3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Footnote row (row 6) - swap the "source" footnote for a
#    "preliminary data" footnote.
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "*алдын алаа маалыматтар"
$ws.Range("B6").Value = "*предварительные данные"
$ws.Range("C6").Value = "*preliminary data"

# ---------------------------------------------------------------------
# 6. Tidy up the saved selection so it points back at A1.
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
